# Yeast_times.xlsx — add a new "Cfinder" row to the community-detection
# table (Tabela1) while keeping the existing "LPANNI" row, pushing it down
# one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1) Duplicate the last data row (13: "LPANNI") into the new row 14,
#    carrying over its formatting/styles exactly.
$ws.Range("B13:L13").Copy($ws.Range("B14:L14"))

# 2) Grow the table (Tabela1) so it covers the new row.
$lo.Resize($ws.Range("B3:L14"))

# 3) The original row 13 becomes the new "Cfinder" algorithm entry; row 14
#    keeps the "LPANNI" value that was copied down in step 1.
$ws.Range("B13").Value = "Cfinder"

# 4) Match the workbook's last active-cell selection.
$ws.Range("B12").Select() | Out-Null
